$d = $word.ActiveDocument

# --- Hunk 1: Sacred Flame - delete "damage increases" sentence, trim trailing space ---
$d.Content.Find.Execute("The spell's damage increases by 1d8 when you reach 5th level (2d8), 11th level (3d8), and 17th level (4d8).", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
